$d = $word.ActiveDocument

$replacements = @(
    @("487÷9=54, 1", "299÷6=49, 5"),
    @("693÷6=115, 3", "459÷6=76, 3"),
    @("149÷4=37, 1", "513÷7=73, 2"),
    @("131÷6=21, 5", "657÷8=82, 1"),
    @("176÷3=58, 2", "327÷5=65, 2"),
    @("483÷3=161, 0", "777÷3=259, 0"),
    @("862÷7=123, 1", "178÷3=59, 1"),
    @("401÷3=133, 2", "762÷4=190, 2"),
    @("603÷9=67, 0", "624÷5=124, 4"),
    @("735÷2=367, 1", "278÷8=34, 6"),
    @("669÷2=334, 1", "365÷4=91, 1"),
    @("555÷7=79, 2", "332÷3=110, 2"),
    @("581÷5=116, 1", "644÷7=92, 0"),
    @("160÷9=17, 7", "636÷9=70, 6"),
    @("658÷3=219, 1", "775÷6=129, 1"),
    @("152÷4=38, 0", "597÷5=119, 2"),
    @("992÷9=110, 2", "557÷2=278, 1"),
    @("939÷9=104, 3", "970÷6=161, 4"),
    @("973÷9=108, 1", "675÷2=337, 1"),
    @("116÷6=19, 2", "894÷8=111, 6"),
    @("108÷8=13, 4", "491÷7=70, 1"),
    @("765÷7=109, 2", "777÷5=155, 2"),
    @("955÷2=477, 1", "744÷7=106, 2"),
    @("205÷2=102, 1", "891÷8=111, 3"),
    @("480÷6=80, 0", "726÷4=181, 2")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
